$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Name of Algo: refresh imputed values (columns A and E) with new results
$ws.Range("E3").Value = 16.15
$ws.Range("E4").Value = 16.3582
$ws.Range("A11").Value = -21.83800000000001
$ws.Range("A12").Value = -21.5329
$ws.Range("E14").Value = 16.4986
$ws.Range("A15").Value = -21.93300000000001
$ws.Range("E26").Value = 16.0421
$ws.Range("A27").Value = -21.86129999999999
$ws.Range("A28").Value = -21.8793
$ws.Range("A31").Value = -21.89500000000002
$ws.Range("E31").Value = 16.6153
$ws.Range("A32").Value = -21.3512
$ws.Range("E35").Value = 16.52
$ws.Range("A36").Value = -19.8694
$ws.Range("E37").Value = 16.59750000000001
$ws.Range("A38").Value = -19.08679999999999
$ws.Range("E39").Value = 16.16069999999999
$ws.Range("E40").Value = 17.01490000000001
$ws.Range("E45").Value = 16.7553
$ws.Range("A46").Value = -21.77290000000001
$ws.Range("E52").Value = 17.09390000000001
$ws.Range("A54").Value = -21.62589999999999
$ws.Range("A55").Value = -22.43880000000001
$ws.Range("A56").Value = -22.1071
$ws.Range("E57").Value = 16.6999
$ws.Range("A67").Value = -21.47319999999998
$ws.Range("A69").Value = -21.71999999999997
$ws.Range("A72").Value = -21.74899999999999
$ws.Range("A73").Value = -19.59529999999998
$ws.Range("E81").Value = 16.49109999999999
$ws.Range("A83").Value = -21.9495
$ws.Range("E83").Value = 16.499
$ws.Range("A86").Value = -21.9605
$ws.Range("A91").Value = -21.41920000000001
$ws.Range("A93").Value = -21.4182
$ws.Range("A99").Value = -19.93849999999999
$ws.Range("E100").Value = 16.25599999999999
$ws.Range("E102").Value = 16.55159999999999
